# Insert a new weekly price-record row right before the current row 300.
# This shifts all existing rows 300:394 down to 301:395 (dimension grows
# from A1:R394 to A1:R395) and leaves a brand-new, empty row 300 which we
# then populate with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(300).Insert()

$ws.Cells.Item(300, 1).Value = 11
$ws.Cells.Item(300, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(300, 3).Value = "Bíobío"
$ws.Cells.Item(300, 4).Value = 44876
$ws.Cells.Item(300, 5).Value = 8
$ws.Cells.Item(300, 6).Value = 100114014
$ws.Cells.Item(300, 7).Value = "Betarraga"
$ws.Cells.Item(300, 8).Value = "Sin especificar"
$ws.Cells.Item(300, 9).Value = "Primera"
$ws.Cells.Item(300, 10).Value = 450
$ws.Cells.Item(300, 11).Value = 700
$ws.Cells.Item(300, 12).Value = 750
$ws.Cells.Item(300, 13).Value = 728
$ws.Cells.Item(300, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(300, 15).Value = "Región Metropolitana"
$ws.Cells.Item(300, 16).Value = 146
$ws.Cells.Item(300, 17).Value = 5
$ws.Cells.Item(300, 18).Value = "Hortaliza"
